$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shifted rows 333-422 (date, calidad, volumen, precio min/max/prom, precio $/Kg)
$ws.Cells.Item(333,4).Value = 44782
$ws.Cells.Item(333,9).Value = 'Segunda'
$ws.Cells.Item(333,10).Value = 600
$ws.Cells.Item(333,11).Value = 1200
$ws.Cells.Item(333,12).Value = 1300
$ws.Cells.Item(333,13).Value = 1250
$ws.Cells.Item(333,16).Value = 1250
$ws.Cells.Item(334,4).Value = 44782
$ws.Cells.Item(334,9).Value = 'Tercera'
$ws.Cells.Item(334,10).Value = 1200
$ws.Cells.Item(334,11).Value = 900
$ws.Cells.Item(334,12).Value = 1000
$ws.Cells.Item(334,13).Value = 950
$ws.Cells.Item(334,16).Value = 950
$ws.Cells.Item(335,4).Value = 44200
$ws.Cells.Item(335,9).Value = 'Segunda'
$ws.Cells.Item(335,10).Value = 1200
$ws.Cells.Item(335,11).Value = 450
$ws.Cells.Item(335,12).Value = 500
$ws.Cells.Item(335,13).Value = 475
$ws.Cells.Item(335,16).Value = 475
$ws.Cells.Item(336,4).Value = 44200
$ws.Cells.Item(336,9).Value = 'Tercera'
$ws.Cells.Item(336,10).Value = 1340
$ws.Cells.Item(336,11).Value = 350
$ws.Cells.Item(336,12).Value = 400
$ws.Cells.Item(336,13).Value = 375
$ws.Cells.Item(336,16).Value = 375
$ws.Cells.Item(337,4).Value = 44748
$ws.Cells.Item(337,9).Value = 'Segunda'
$ws.Cells.Item(337,10).Value = 1200
$ws.Cells.Item(337,11).Value = 600
$ws.Cells.Item(337,12).Value = 700
$ws.Cells.Item(337,13).Value = 650
$ws.Cells.Item(337,16).Value = 650
$ws.Cells.Item(338,4).Value = 44748
$ws.Cells.Item(338,9).Value = 'Tercera'
$ws.Cells.Item(338,10).Value = 1200
$ws.Cells.Item(338,11).Value = 400
$ws.Cells.Item(338,12).Value = 500
$ws.Cells.Item(338,13).Value = 450
$ws.Cells.Item(338,16).Value = 450
$ws.Cells.Item(339,4).Value = 44424
$ws.Cells.Item(339,9).Value = 'Segunda'
$ws.Cells.Item(339,10).Value = 1000
$ws.Cells.Item(339,11).Value = 800
$ws.Cells.Item(339,12).Value = 900
$ws.Cells.Item(339,13).Value = 850
$ws.Cells.Item(339,16).Value = 850
$ws.Cells.Item(340,4).Value = 44424
$ws.Cells.Item(340,9).Value = 'Tercera'
$ws.Cells.Item(340,10).Value = 1200
$ws.Cells.Item(340,11).Value = 500
$ws.Cells.Item(340,12).Value = 600
$ws.Cells.Item(340,13).Value = 550
$ws.Cells.Item(340,16).Value = 550
$ws.Cells.Item(341,4).Value = 44259
$ws.Cells.Item(341,9).Value = 'Tercera'
$ws.Cells.Item(341,10).Value = 1200
$ws.Cells.Item(341,11).Value = 700
$ws.Cells.Item(341,12).Value = 800
$ws.Cells.Item(341,13).Value = 750
$ws.Cells.Item(341,16).Value = 750
$ws.Cells.Item(342,4).Value = 44202
$ws.Cells.Item(342,9).Value = 'Tercera'
$ws.Cells.Item(342,10).Value = 1500
$ws.Cells.Item(342,11).Value = 300
$ws.Cells.Item(342,12).Value = 350
$ws.Cells.Item(342,13).Value = 325
$ws.Cells.Item(342,16).Value = 325
$ws.Cells.Item(343,4).Value = 44371
$ws.Cells.Item(343,9).Value = 'Tercera'
$ws.Cells.Item(343,10).Value = 1000
$ws.Cells.Item(343,11).Value = 700
$ws.Cells.Item(343,12).Value = 800
$ws.Cells.Item(343,13).Value = 760
$ws.Cells.Item(343,16).Value = 760
$ws.Cells.Item(344,4).Value = 44162
$ws.Cells.Item(344,9).Value = 'Segunda'
$ws.Cells.Item(344,10).Value = 1200
$ws.Cells.Item(344,11).Value = 450
$ws.Cells.Item(344,12).Value = 500
$ws.Cells.Item(344,13).Value = 475
$ws.Cells.Item(344,16).Value = 475
$ws.Cells.Item(345,4).Value = 44162
$ws.Cells.Item(345,9).Value = 'Tercera'
$ws.Cells.Item(345,10).Value = 1500
$ws.Cells.Item(345,11).Value = 350
$ws.Cells.Item(345,12).Value = 400
$ws.Cells.Item(345,13).Value = 375
$ws.Cells.Item(345,16).Value = 375
$ws.Cells.Item(346,4).Value = 44235
$ws.Cells.Item(346,9).Value = 'Segunda'
$ws.Cells.Item(346,10).Value = 700
$ws.Cells.Item(346,11).Value = 600
$ws.Cells.Item(346,12).Value = 700
$ws.Cells.Item(346,13).Value = 650
$ws.Cells.Item(346,16).Value = 650
$ws.Cells.Item(347,4).Value = 44235
$ws.Cells.Item(347,9).Value = 'Tercera'
$ws.Cells.Item(347,10).Value = 1200
$ws.Cells.Item(347,11).Value = 450
$ws.Cells.Item(347,12).Value = 500
$ws.Cells.Item(347,13).Value = 475
$ws.Cells.Item(347,16).Value = 475
$ws.Cells.Item(348,4).Value = 44726
$ws.Cells.Item(348,9).Value = 'Segunda'
$ws.Cells.Item(348,10).Value = 1200
$ws.Cells.Item(348,11).Value = 500
$ws.Cells.Item(348,12).Value = 600
$ws.Cells.Item(348,13).Value = 550
$ws.Cells.Item(348,16).Value = 550
$ws.Cells.Item(349,4).Value = 44726
$ws.Cells.Item(349,9).Value = 'Tercera'
$ws.Cells.Item(349,10).Value = 1300
$ws.Cells.Item(349,11).Value = 400
$ws.Cells.Item(349,12).Value = 450
$ws.Cells.Item(349,13).Value = 425
$ws.Cells.Item(349,16).Value = 425
$ws.Cells.Item(350,4).Value = 44427
$ws.Cells.Item(350,9).Value = 'Segunda'
$ws.Cells.Item(350,10).Value = 1300
$ws.Cells.Item(350,11).Value = 700
$ws.Cells.Item(350,12).Value = 800
$ws.Cells.Item(350,13).Value = 750
$ws.Cells.Item(350,16).Value = 750
$ws.Cells.Item(351,4).Value = 44427
$ws.Cells.Item(351,9).Value = 'Tercera'
$ws.Cells.Item(351,10).Value = 1200
$ws.Cells.Item(351,11).Value = 500
$ws.Cells.Item(351,12).Value = 600
$ws.Cells.Item(351,13).Value = 550
$ws.Cells.Item(351,16).Value = 550
$ws.Cells.Item(352,4).Value = 44441
$ws.Cells.Item(352,9).Value = 'Segunda'
$ws.Cells.Item(352,10).Value = 800
$ws.Cells.Item(352,11).Value = 600
$ws.Cells.Item(352,12).Value = 700
$ws.Cells.Item(352,13).Value = 650
$ws.Cells.Item(352,16).Value = 650
$ws.Cells.Item(353,4).Value = 44441
$ws.Cells.Item(353,9).Value = 'Tercera'
$ws.Cells.Item(353,10).Value = 700
$ws.Cells.Item(353,11).Value = 500
$ws.Cells.Item(353,12).Value = 600
$ws.Cells.Item(353,13).Value = 550
$ws.Cells.Item(353,16).Value = 550
$ws.Cells.Item(354,4).Value = 44174
$ws.Cells.Item(354,9).Value = 'Primera'
$ws.Cells.Item(354,10).Value = 500
$ws.Cells.Item(354,11).Value = 450
$ws.Cells.Item(354,12).Value = 550
$ws.Cells.Item(354,13).Value = 500
$ws.Cells.Item(354,16).Value = 500
$ws.Cells.Item(355,4).Value = 44174
$ws.Cells.Item(355,9).Value = 'Segunda'
$ws.Cells.Item(355,10).Value = 1200
$ws.Cells.Item(355,11).Value = 400
$ws.Cells.Item(355,12).Value = 450
$ws.Cells.Item(355,13).Value = 425
$ws.Cells.Item(355,16).Value = 425
$ws.Cells.Item(356,4).Value = 44174
$ws.Cells.Item(356,9).Value = 'Tercera'
$ws.Cells.Item(356,10).Value = 1200
$ws.Cells.Item(356,11).Value = 300
$ws.Cells.Item(356,12).Value = 350
$ws.Cells.Item(356,13).Value = 325
$ws.Cells.Item(356,16).Value = 325
$ws.Cells.Item(357,4).Value = 44419
$ws.Cells.Item(357,9).Value = 'Segunda'
$ws.Cells.Item(357,10).Value = 500
$ws.Cells.Item(357,11).Value = 700
$ws.Cells.Item(357,12).Value = 800
$ws.Cells.Item(357,13).Value = 750
$ws.Cells.Item(357,16).Value = 750
$ws.Cells.Item(358,4).Value = 44419
$ws.Cells.Item(358,9).Value = 'Tercera'
$ws.Cells.Item(358,10).Value = 1200
$ws.Cells.Item(358,11).Value = 400
$ws.Cells.Item(358,12).Value = 500
$ws.Cells.Item(358,13).Value = 450
$ws.Cells.Item(358,16).Value = 450
$ws.Cells.Item(359,4).Value = 44280
$ws.Cells.Item(359,9).Value = 'Segunda'
$ws.Cells.Item(359,10).Value = 500
$ws.Cells.Item(359,11).Value = 900
$ws.Cells.Item(359,12).Value = 1000
$ws.Cells.Item(359,13).Value = 950
$ws.Cells.Item(359,16).Value = 950
$ws.Cells.Item(360,4).Value = 44280
$ws.Cells.Item(360,9).Value = 'Tercera'
$ws.Cells.Item(360,10).Value = 700
$ws.Cells.Item(360,11).Value = 800
$ws.Cells.Item(360,12).Value = 900
$ws.Cells.Item(360,13).Value = 850
$ws.Cells.Item(360,16).Value = 850
$ws.Cells.Item(361,4).Value = 44412
$ws.Cells.Item(361,9).Value = 'Primera'
$ws.Cells.Item(361,10).Value = 500
$ws.Cells.Item(361,11).Value = 900
$ws.Cells.Item(361,12).Value = 1000
$ws.Cells.Item(361,13).Value = 950
$ws.Cells.Item(361,16).Value = 950
$ws.Cells.Item(362,4).Value = 44412
$ws.Cells.Item(362,9).Value = 'Segunda'
$ws.Cells.Item(362,10).Value = 1200
$ws.Cells.Item(362,11).Value = 800
$ws.Cells.Item(362,12).Value = 900
$ws.Cells.Item(362,13).Value = 850
$ws.Cells.Item(362,16).Value = 850
$ws.Cells.Item(363,4).Value = 44412
$ws.Cells.Item(363,9).Value = 'Tercera'
$ws.Cells.Item(363,10).Value = 1300
$ws.Cells.Item(363,11).Value = 450
$ws.Cells.Item(363,12).Value = 500
$ws.Cells.Item(363,13).Value = 475
$ws.Cells.Item(363,16).Value = 475
$ws.Cells.Item(364,4).Value = 44237
$ws.Cells.Item(364,9).Value = 'Segunda'
$ws.Cells.Item(364,10).Value = 740
$ws.Cells.Item(364,11).Value = 1000
$ws.Cells.Item(364,12).Value = 1200
$ws.Cells.Item(364,13).Value = 1100
$ws.Cells.Item(364,16).Value = 1100
$ws.Cells.Item(365,4).Value = 44237
$ws.Cells.Item(365,9).Value = 'Tercera'
$ws.Cells.Item(365,10).Value = 600
$ws.Cells.Item(365,11).Value = 700
$ws.Cells.Item(365,12).Value = 800
$ws.Cells.Item(365,13).Value = 750
$ws.Cells.Item(365,16).Value = 750
$ws.Cells.Item(366,4).Value = 44628
$ws.Cells.Item(366,9).Value = 'Segunda'
$ws.Cells.Item(366,10).Value = 1200
$ws.Cells.Item(366,11).Value = 450
$ws.Cells.Item(366,12).Value = 500
$ws.Cells.Item(366,13).Value = 475
$ws.Cells.Item(366,16).Value = 475
$ws.Cells.Item(367,4).Value = 44628
$ws.Cells.Item(367,9).Value = 'Tercera'
$ws.Cells.Item(367,10).Value = 1200
$ws.Cells.Item(367,11).Value = 350
$ws.Cells.Item(367,12).Value = 400
$ws.Cells.Item(367,13).Value = 375
$ws.Cells.Item(367,16).Value = 375
$ws.Cells.Item(368,4).Value = 44483
$ws.Cells.Item(368,9).Value = 'Primera'
$ws.Cells.Item(368,10).Value = 1000
$ws.Cells.Item(368,11).Value = 500
$ws.Cells.Item(368,12).Value = 600
$ws.Cells.Item(368,13).Value = 550
$ws.Cells.Item(368,16).Value = 550
$ws.Cells.Item(369,4).Value = 44483
$ws.Cells.Item(369,9).Value = 'Segunda'
$ws.Cells.Item(369,10).Value = 1400
$ws.Cells.Item(369,11).Value = 450
$ws.Cells.Item(369,12).Value = 500
$ws.Cells.Item(369,13).Value = 475
$ws.Cells.Item(369,16).Value = 475
$ws.Cells.Item(370,4).Value = 44483
$ws.Cells.Item(370,9).Value = 'Tercera'
$ws.Cells.Item(370,10).Value = 1200
$ws.Cells.Item(370,11).Value = 350
$ws.Cells.Item(370,12).Value = 400
$ws.Cells.Item(370,13).Value = 375
$ws.Cells.Item(370,16).Value = 375
$ws.Cells.Item(371,4).Value = 44175
$ws.Cells.Item(371,9).Value = 'Segunda'
$ws.Cells.Item(371,10).Value = 1200
$ws.Cells.Item(371,11).Value = 350
$ws.Cells.Item(371,12).Value = 400
$ws.Cells.Item(371,13).Value = 375
$ws.Cells.Item(371,16).Value = 375
$ws.Cells.Item(372,4).Value = 44175
$ws.Cells.Item(372,9).Value = 'Tercera'
$ws.Cells.Item(372,10).Value = 1600
$ws.Cells.Item(372,11).Value = 250
$ws.Cells.Item(372,12).Value = 300
$ws.Cells.Item(372,13).Value = 275
$ws.Cells.Item(372,16).Value = 275
$ws.Cells.Item(373,4).Value = 44469
$ws.Cells.Item(373,9).Value = 'Segunda'
$ws.Cells.Item(373,10).Value = 700
$ws.Cells.Item(373,11).Value = 600
$ws.Cells.Item(373,12).Value = 700
$ws.Cells.Item(373,13).Value = 650
$ws.Cells.Item(373,16).Value = 650
$ws.Cells.Item(374,4).Value = 44469
$ws.Cells.Item(374,9).Value = 'Tercera'
$ws.Cells.Item(374,10).Value = 1200
$ws.Cells.Item(374,11).Value = 400
$ws.Cells.Item(374,12).Value = 500
$ws.Cells.Item(374,13).Value = 450
$ws.Cells.Item(374,16).Value = 450
$ws.Cells.Item(375,4).Value = 44434
$ws.Cells.Item(375,9).Value = 'Segunda'
$ws.Cells.Item(375,10).Value = 1000
$ws.Cells.Item(375,11).Value = 900
$ws.Cells.Item(375,12).Value = 1000
$ws.Cells.Item(375,13).Value = 950
$ws.Cells.Item(375,16).Value = 950
$ws.Cells.Item(376,4).Value = 44434
$ws.Cells.Item(376,9).Value = 'Tercera'
$ws.Cells.Item(376,10).Value = 900
$ws.Cells.Item(376,11).Value = 600
$ws.Cells.Item(376,12).Value = 700
$ws.Cells.Item(376,13).Value = 650
$ws.Cells.Item(376,16).Value = 650
$ws.Cells.Item(377,4).Value = 44776
$ws.Cells.Item(377,9).Value = 'Tercera'
$ws.Cells.Item(377,10).Value = 800
$ws.Cells.Item(377,11).Value = 450
$ws.Cells.Item(377,12).Value = 500
$ws.Cells.Item(377,13).Value = 475
$ws.Cells.Item(377,16).Value = 475
$ws.Cells.Item(378,4).Value = 44494
$ws.Cells.Item(378,9).Value = 'Segunda'
$ws.Cells.Item(378,10).Value = 1200
$ws.Cells.Item(378,11).Value = 400
$ws.Cells.Item(378,12).Value = 500
$ws.Cells.Item(378,13).Value = 450
$ws.Cells.Item(378,16).Value = 450
$ws.Cells.Item(379,4).Value = 44494
$ws.Cells.Item(379,9).Value = 'Tercera'
$ws.Cells.Item(379,10).Value = 800
$ws.Cells.Item(379,11).Value = 300
$ws.Cells.Item(379,12).Value = 350
$ws.Cells.Item(379,13).Value = 325
$ws.Cells.Item(379,16).Value = 325
$ws.Cells.Item(380,4).Value = 44487
$ws.Cells.Item(380,9).Value = 'Primera'
$ws.Cells.Item(380,10).Value = 700
$ws.Cells.Item(380,11).Value = 600
$ws.Cells.Item(380,12).Value = 700
$ws.Cells.Item(380,13).Value = 650
$ws.Cells.Item(380,16).Value = 650
$ws.Cells.Item(381,4).Value = 44487
$ws.Cells.Item(381,9).Value = 'Segunda'
$ws.Cells.Item(381,10).Value = 1200
$ws.Cells.Item(381,11).Value = 500
$ws.Cells.Item(381,12).Value = 550
$ws.Cells.Item(381,13).Value = 525
$ws.Cells.Item(381,16).Value = 525
$ws.Cells.Item(382,4).Value = 44487
$ws.Cells.Item(382,9).Value = 'Tercera'
$ws.Cells.Item(382,10).Value = 1200
$ws.Cells.Item(382,11).Value = 350
$ws.Cells.Item(382,12).Value = 400
$ws.Cells.Item(382,13).Value = 375
$ws.Cells.Item(382,16).Value = 375
$ws.Cells.Item(383,4).Value = 44266
$ws.Cells.Item(383,9).Value = 'Segunda'
$ws.Cells.Item(383,10).Value = 800
$ws.Cells.Item(383,11).Value = 800
$ws.Cells.Item(383,12).Value = 850
$ws.Cells.Item(383,13).Value = 825
$ws.Cells.Item(383,16).Value = 825
$ws.Cells.Item(384,4).Value = 44266
$ws.Cells.Item(384,9).Value = 'Tercera'
$ws.Cells.Item(384,10).Value = 900
$ws.Cells.Item(384,11).Value = 700
$ws.Cells.Item(384,12).Value = 750
$ws.Cells.Item(384,13).Value = 725
$ws.Cells.Item(384,16).Value = 725
$ws.Cells.Item(385,4).Value = 44488
$ws.Cells.Item(385,9).Value = 'Primera'
$ws.Cells.Item(385,10).Value = 1000
$ws.Cells.Item(385,11).Value = 600
$ws.Cells.Item(385,12).Value = 700
$ws.Cells.Item(385,13).Value = 650
$ws.Cells.Item(385,16).Value = 650
$ws.Cells.Item(386,4).Value = 44488
$ws.Cells.Item(386,9).Value = 'Segunda'
$ws.Cells.Item(386,10).Value = 1200
$ws.Cells.Item(386,11).Value = 400
$ws.Cells.Item(386,12).Value = 500
$ws.Cells.Item(386,13).Value = 450
$ws.Cells.Item(386,16).Value = 450
$ws.Cells.Item(387,4).Value = 44488
$ws.Cells.Item(387,9).Value = 'Tercera'
$ws.Cells.Item(387,10).Value = 1200
$ws.Cells.Item(387,11).Value = 300
$ws.Cells.Item(387,12).Value = 350
$ws.Cells.Item(387,13).Value = 325
$ws.Cells.Item(387,16).Value = 325
$ws.Cells.Item(388,4).Value = 44238
$ws.Cells.Item(388,9).Value = 'Tercera'
$ws.Cells.Item(388,10).Value = 700
$ws.Cells.Item(388,11).Value = 900
$ws.Cells.Item(388,12).Value = 1000
$ws.Cells.Item(388,13).Value = 950
$ws.Cells.Item(388,16).Value = 950
$ws.Cells.Item(389,4).Value = 44340
$ws.Cells.Item(389,9).Value = 'Segunda'
$ws.Cells.Item(389,10).Value = 700
$ws.Cells.Item(389,11).Value = 850
$ws.Cells.Item(389,12).Value = 950
$ws.Cells.Item(389,13).Value = 900
$ws.Cells.Item(389,16).Value = 900
$ws.Cells.Item(390,4).Value = 44340
$ws.Cells.Item(390,9).Value = 'Tercera'
$ws.Cells.Item(390,10).Value = 1200
$ws.Cells.Item(390,11).Value = 350
$ws.Cells.Item(390,12).Value = 450
$ws.Cells.Item(390,13).Value = 400
$ws.Cells.Item(390,16).Value = 400
$ws.Cells.Item(391,4).Value = 44496
$ws.Cells.Item(391,9).Value = 'Primera'
$ws.Cells.Item(391,10).Value = 1000
$ws.Cells.Item(391,11).Value = 450
$ws.Cells.Item(391,12).Value = 500
$ws.Cells.Item(391,13).Value = 475
$ws.Cells.Item(391,16).Value = 475
$ws.Cells.Item(392,4).Value = 44496
$ws.Cells.Item(392,9).Value = 'Segunda'
$ws.Cells.Item(392,10).Value = 1500
$ws.Cells.Item(392,11).Value = 350
$ws.Cells.Item(392,12).Value = 450
$ws.Cells.Item(392,13).Value = 400
$ws.Cells.Item(392,16).Value = 400
$ws.Cells.Item(393,4).Value = 44399
$ws.Cells.Item(393,9).Value = 'Segunda'
$ws.Cells.Item(393,10).Value = 800
$ws.Cells.Item(393,11).Value = 1000
$ws.Cells.Item(393,12).Value = 1200
$ws.Cells.Item(393,13).Value = 1100
$ws.Cells.Item(393,16).Value = 1100
$ws.Cells.Item(394,4).Value = 44399
$ws.Cells.Item(394,9).Value = 'Tercera'
$ws.Cells.Item(394,10).Value = 1200
$ws.Cells.Item(394,11).Value = 500
$ws.Cells.Item(394,12).Value = 600
$ws.Cells.Item(394,13).Value = 550
$ws.Cells.Item(394,16).Value = 550
$ws.Cells.Item(395,4).Value = 44425
$ws.Cells.Item(395,9).Value = 'Segunda'
$ws.Cells.Item(395,10).Value = 800
$ws.Cells.Item(395,11).Value = 750
$ws.Cells.Item(395,12).Value = 800
$ws.Cells.Item(395,13).Value = 775
$ws.Cells.Item(395,16).Value = 775
$ws.Cells.Item(396,4).Value = 44425
$ws.Cells.Item(396,9).Value = 'Tercera'
$ws.Cells.Item(396,10).Value = 1200
$ws.Cells.Item(396,11).Value = 450
$ws.Cells.Item(396,12).Value = 500
$ws.Cells.Item(396,13).Value = 475
$ws.Cells.Item(396,16).Value = 475
$ws.Cells.Item(397,4).Value = 44377
$ws.Cells.Item(397,9).Value = 'Tercera'
$ws.Cells.Item(397,10).Value = 850
$ws.Cells.Item(397,11).Value = 500
$ws.Cells.Item(397,12).Value = 600
$ws.Cells.Item(397,13).Value = 547
$ws.Cells.Item(397,16).Value = 547
$ws.Cells.Item(398,4).Value = 44181
$ws.Cells.Item(398,9).Value = 'Segunda'
$ws.Cells.Item(398,10).Value = 960
$ws.Cells.Item(398,11).Value = 300
$ws.Cells.Item(398,12).Value = 350
$ws.Cells.Item(398,13).Value = 325
$ws.Cells.Item(398,16).Value = 325
$ws.Cells.Item(399,4).Value = 44181
$ws.Cells.Item(399,9).Value = 'Tercera'
$ws.Cells.Item(399,10).Value = 1200
$ws.Cells.Item(399,11).Value = 250
$ws.Cells.Item(399,12).Value = 300
$ws.Cells.Item(399,13).Value = 275
$ws.Cells.Item(399,16).Value = 275
$ws.Cells.Item(400,4).Value = 44497
$ws.Cells.Item(400,9).Value = 'Segunda'
$ws.Cells.Item(400,10).Value = 1200
$ws.Cells.Item(400,11).Value = 400
$ws.Cells.Item(400,12).Value = 500
$ws.Cells.Item(400,13).Value = 450
$ws.Cells.Item(400,16).Value = 450
$ws.Cells.Item(401,4).Value = 44497
$ws.Cells.Item(401,9).Value = 'Tercera'
$ws.Cells.Item(401,10).Value = 1200
$ws.Cells.Item(401,11).Value = 350
$ws.Cells.Item(401,12).Value = 400
$ws.Cells.Item(401,13).Value = 375
$ws.Cells.Item(401,16).Value = 375
$ws.Cells.Item(402,4).Value = 44362
$ws.Cells.Item(402,9).Value = 'Segunda'
$ws.Cells.Item(402,10).Value = 1000
$ws.Cells.Item(402,11).Value = 500
$ws.Cells.Item(402,12).Value = 600
$ws.Cells.Item(402,13).Value = 550
$ws.Cells.Item(402,16).Value = 550
$ws.Cells.Item(403,4).Value = 44362
$ws.Cells.Item(403,9).Value = 'Tercera'
$ws.Cells.Item(403,10).Value = 1000
$ws.Cells.Item(403,11).Value = 400
$ws.Cells.Item(403,12).Value = 500
$ws.Cells.Item(403,13).Value = 450
$ws.Cells.Item(403,16).Value = 450
$ws.Cells.Item(404,4).Value = 44557
$ws.Cells.Item(404,9).Value = 'Segunda'
$ws.Cells.Item(404,10).Value = 1300
$ws.Cells.Item(404,11).Value = 350
$ws.Cells.Item(404,12).Value = 400
$ws.Cells.Item(404,13).Value = 375
$ws.Cells.Item(404,16).Value = 375
$ws.Cells.Item(405,4).Value = 44557
$ws.Cells.Item(405,9).Value = 'Tercera'
$ws.Cells.Item(405,10).Value = 1200
$ws.Cells.Item(405,11).Value = 300
$ws.Cells.Item(405,12).Value = 350
$ws.Cells.Item(405,13).Value = 325
$ws.Cells.Item(405,16).Value = 325
$ws.Cells.Item(406,4).Value = 44747
$ws.Cells.Item(406,9).Value = 'Segunda'
$ws.Cells.Item(406,10).Value = 1200
$ws.Cells.Item(406,11).Value = 500
$ws.Cells.Item(406,12).Value = 600
$ws.Cells.Item(406,13).Value = 550
$ws.Cells.Item(406,16).Value = 550
$ws.Cells.Item(407,4).Value = 44747
$ws.Cells.Item(407,9).Value = 'Tercera'
$ws.Cells.Item(407,10).Value = 1200
$ws.Cells.Item(407,11).Value = 450
$ws.Cells.Item(407,12).Value = 500
$ws.Cells.Item(407,13).Value = 475
$ws.Cells.Item(407,16).Value = 475
$ws.Cells.Item(408,4).Value = 44357
$ws.Cells.Item(408,9).Value = 'Segunda'
$ws.Cells.Item(408,10).Value = 700
$ws.Cells.Item(408,11).Value = 800
$ws.Cells.Item(408,12).Value = 900
$ws.Cells.Item(408,13).Value = 850
$ws.Cells.Item(408,16).Value = 850
$ws.Cells.Item(409,4).Value = 44357
$ws.Cells.Item(409,9).Value = 'Tercera'
$ws.Cells.Item(409,10).Value = 600
$ws.Cells.Item(409,11).Value = 600
$ws.Cells.Item(409,12).Value = 700
$ws.Cells.Item(409,13).Value = 650
$ws.Cells.Item(409,16).Value = 650
$ws.Cells.Item(410,4).Value = 44551
$ws.Cells.Item(410,9).Value = 'Primera'
$ws.Cells.Item(410,10).Value = 1000
$ws.Cells.Item(410,11).Value = 400
$ws.Cells.Item(410,12).Value = 500
$ws.Cells.Item(410,13).Value = 450
$ws.Cells.Item(410,16).Value = 450
$ws.Cells.Item(411,4).Value = 44551
$ws.Cells.Item(411,9).Value = 'Segunda'
$ws.Cells.Item(411,10).Value = 1200
$ws.Cells.Item(411,11).Value = 400
$ws.Cells.Item(411,12).Value = 450
$ws.Cells.Item(411,13).Value = 425
$ws.Cells.Item(411,16).Value = 425
$ws.Cells.Item(412,4).Value = 44551
$ws.Cells.Item(412,9).Value = 'Tercera'
$ws.Cells.Item(412,10).Value = 1000
$ws.Cells.Item(412,11).Value = 300
$ws.Cells.Item(412,12).Value = 350
$ws.Cells.Item(412,13).Value = 325
$ws.Cells.Item(412,16).Value = 325
$ws.Cells.Item(413,4).Value = 44355
$ws.Cells.Item(413,9).Value = 'Segunda'
$ws.Cells.Item(413,10).Value = 1000
$ws.Cells.Item(413,11).Value = 600
$ws.Cells.Item(413,12).Value = 700
$ws.Cells.Item(413,13).Value = 650
$ws.Cells.Item(413,16).Value = 650
$ws.Cells.Item(414,4).Value = 44355
$ws.Cells.Item(414,9).Value = 'Tercera'
$ws.Cells.Item(414,10).Value = 1000
$ws.Cells.Item(414,11).Value = 400
$ws.Cells.Item(414,12).Value = 500
$ws.Cells.Item(414,13).Value = 450
$ws.Cells.Item(414,16).Value = 450
$ws.Cells.Item(415,4).Value = 44391
$ws.Cells.Item(415,9).Value = 'Segunda'
$ws.Cells.Item(415,10).Value = 900
$ws.Cells.Item(415,11).Value = 900
$ws.Cells.Item(415,12).Value = 1000
$ws.Cells.Item(415,13).Value = 950
$ws.Cells.Item(415,16).Value = 950
$ws.Cells.Item(416,4).Value = 44391
$ws.Cells.Item(416,9).Value = 'Tercera'
$ws.Cells.Item(416,10).Value = 1000
$ws.Cells.Item(416,11).Value = 600
$ws.Cells.Item(416,12).Value = 700
$ws.Cells.Item(416,13).Value = 650
$ws.Cells.Item(416,16).Value = 650
$ws.Cells.Item(417,4).Value = 44453
$ws.Cells.Item(417,9).Value = 'Segunda'
$ws.Cells.Item(417,10).Value = 1200
$ws.Cells.Item(417,11).Value = 700
$ws.Cells.Item(417,12).Value = 800
$ws.Cells.Item(417,13).Value = 750
$ws.Cells.Item(417,16).Value = 750
$ws.Cells.Item(418,4).Value = 44453
$ws.Cells.Item(418,9).Value = 'Tercera'
$ws.Cells.Item(418,10).Value = 1300
$ws.Cells.Item(418,11).Value = 400
$ws.Cells.Item(418,12).Value = 500
$ws.Cells.Item(418,13).Value = 450
$ws.Cells.Item(418,16).Value = 450
$ws.Cells.Item(419,4).Value = 44609
$ws.Cells.Item(419,9).Value = 'Tercera'
$ws.Cells.Item(419,10).Value = 1600
$ws.Cells.Item(419,11).Value = 350
$ws.Cells.Item(419,12).Value = 400
$ws.Cells.Item(419,13).Value = 375
$ws.Cells.Item(419,16).Value = 375
$ws.Cells.Item(420,4).Value = 44489
$ws.Cells.Item(420,9).Value = 'Primera'
$ws.Cells.Item(420,10).Value = 900
$ws.Cells.Item(420,11).Value = 500
$ws.Cells.Item(420,12).Value = 600
$ws.Cells.Item(420,13).Value = 550
$ws.Cells.Item(420,16).Value = 550
$ws.Cells.Item(421,4).Value = 44489
$ws.Cells.Item(421,9).Value = 'Segunda'
$ws.Cells.Item(421,10).Value = 1200
$ws.Cells.Item(421,11).Value = 400
$ws.Cells.Item(421,12).Value = 500
$ws.Cells.Item(421,13).Value = 450
$ws.Cells.Item(421,16).Value = 450
$ws.Cells.Item(422,4).Value = 44489
$ws.Cells.Item(422,9).Value = 'Tercera'
$ws.Cells.Item(422,10).Value = 1200
$ws.Cells.Item(422,11).Value = 350
$ws.Cells.Item(422,12).Value = 400
$ws.Cells.Item(422,13).Value = 375
$ws.Cells.Item(422,16).Value = 375

# New rows 423-424 (full rows, including constant columns)
$ws.Cells.Item(423,1).Value = 1
$ws.Cells.Item(423,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(423,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(423,5).Value = 15
$ws.Cells.Item(423,6).Value = 100112023
$ws.Cells.Item(423,7).Value = 'Brócoli'
$ws.Cells.Item(423,8).Value = 'Sin especificar'
$ws.Cells.Item(423,14).Value = '$/unidad'
$ws.Cells.Item(423,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(423,17).Value = 1
$ws.Cells.Item(423,18).Value = 'Hortaliza'
$ws.Cells.Item(423,4).Value = 44572
$ws.Cells.Item(423,9).Value = 'Segunda'
$ws.Cells.Item(423,10).Value = 1400
$ws.Cells.Item(423,11).Value = 350
$ws.Cells.Item(423,12).Value = 400
$ws.Cells.Item(423,13).Value = 375
$ws.Cells.Item(423,16).Value = 375
$ws.Cells.Item(423,4).NumberFormat = $ws.Cells.Item(421,4).NumberFormat
$ws.Cells.Item(424,1).Value = 1
$ws.Cells.Item(424,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(424,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(424,5).Value = 15
$ws.Cells.Item(424,6).Value = 100112023
$ws.Cells.Item(424,7).Value = 'Brócoli'
$ws.Cells.Item(424,8).Value = 'Sin especificar'
$ws.Cells.Item(424,14).Value = '$/unidad'
$ws.Cells.Item(424,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(424,17).Value = 1
$ws.Cells.Item(424,18).Value = 'Hortaliza'
$ws.Cells.Item(424,4).Value = 44572
$ws.Cells.Item(424,9).Value = 'Tercera'
$ws.Cells.Item(424,10).Value = 1300
$ws.Cells.Item(424,11).Value = 300
$ws.Cells.Item(424,12).Value = 350
$ws.Cells.Item(424,13).Value = 325
$ws.Cells.Item(424,16).Value = 325
$ws.Cells.Item(424,4).NumberFormat = $ws.Cells.Item(422,4).NumberFormat
